# Ajuste Registro de Cotizacion Elementos Relacionados Bien/Servicio
#
# The "Servicio" sheet used to carry extra columns (E:H) for capturing a
# unit-of-measure + "Tiempo de Ejecucion" (years/months/days) breakdown for
# services, plus a trailing amount column (H). Those extra columns are no
# longer needed; the sheet is trimmed back down to Tipo/Nombre/Descripcion +
# a single "Cantidad" (amount) column, reusing what used to be column H.

$wb = $excel.ActiveWorkbook

# --- Sheet "Servicio": drop the now-unused D:G columns -------------------
# Deleting entire columns D through G shifts the old column H (amount,
# "Campo para cantidad del Servicio" / "Cantidad" / 456465.45) left into the
# new column D, drops the E2:G2 merge + the obsolete "Tiempo Ejecucion"
# validations automatically, and recomputes the sheet dimension.
$wsServicio = $wb.Worksheets.Item("Servicio")
[void]$wsServicio.Range("D1:G4").EntireColumn.Delete()

# Selection moved to D14 in the saved file.
[void]$wsServicio.Range("D14").Select()

# --- Sheet "Bien": selection moved up to A3 -------------------------------
$wsBien = $wb.Worksheets.Item("Bien")
[void]$wsBien.Range("A3").Select()

# Re-activate the originally active sheet ("Bien" is tabSelected="true").
[void]$wsBien.Activate()
